$wb = $excel.ActiveWorkbook

# --- Sheet reference ---
$wsBOM = $wb.Worksheets.Item("BOM Overview")

# --- Insert a new blank row at position 10 (pushes old row10+ down by one) ---
$wsBOM.Rows.Item(10).Insert()

# --- Populate row 9 (previously blank spacer row) with the new line item ---
$wsBOM.Range("A9").Value = "Cables and Connectors Estimated Price"
$wsBOM.Range("B9").Value = 50

# --- Nudge the floating shapes anchored below the insertion point down by
#     one row height so their anchors follow the shifted rows (the engine's
#     row-insert does not itself relocate drawing anchors) ---
$wsBOM.Shapes.Item("Picture 1").Top = 245.1
$wsBOM.Shapes.Item("TextBox 3").Top = 215.85

# --- Activate "BOM Overview" sheet (was "Complete Part List") ---
$wsBOM.Activate()

